$wb = $excel.ActiveWorkbook

# --- Resumen sheet: the winning zone changed from Z2 to Z1 (Maximo stays 361,0) ---
$wsResumen = $wb.Worksheets.Item("Resumen")
$wsResumen.Range("B2").Value = "Z1"
$wsResumen.Range("C2").Value = "361,0"

# --- Solucion sheet: every "Salida" assignment shifts by 20 positions (mod 40) ---
$wsSolucion = $wb.Worksheets.Item("Solucion")
$salidas = @("S021","S001","S011","S031","S002","S022","S012","S032","S003","S023","S013","S033","S004","S024","S014","S034","S025","S005","S015","S035","S006","S026","S016","S036","S027","S007","S017","S037","S028","S008","S018","S038","S029","S009","S039","S019","S010","S030","S040","S020")
for ($i = 0; $i -lt $salidas.Length; $i++) {
    $row = $i + 2
    $wsSolucion.Cells.Item($row, 2).Value = $salidas[$i]
}

# --- Metricas sheet: Z1 and Z2 swap their Tiempo values ---
$wsMetricas = $wb.Worksheets.Item("Metricas")
$wsMetricas.Range("A2").Value = "Z1"
$wsMetricas.Range("B2").Value = "361,0"
$wsMetricas.Range("A3").Value = "Z2"
$wsMetricas.Range("B3").Value = "355,7"
